# Add 5 more rows (rows 4-8) to Sheet1, duplicating the data pattern of row 3,
# as part of building an auto-evaluation pipeline dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$idValue    = "12306-1"
$appValue   = 12306
$queryValue = "订一张上午11点左右从郑州到北京的高铁票，要求二等座F座"

for ($r = 4; $r -le 8; $r++) {
    $ws.Cells.Item($r, 1).Value = $idValue
    $ws.Cells.Item($r, 2).Value = $appValue
    $ws.Cells.Item($r, 3).Value = $queryValue
}

# Match the resulting selection state recorded in the saved workbook.
$ws.Range("A8:C8").Select()
